$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the new "2022-Q1" holdings-detail sheet.
#    The existing "2021-Q4" sheet already has exactly the column layout /
#    styles we need (基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名), so clone it and just overwrite the data rows.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy([System.Type]::Missing, $q4)
$newSheet = $wb.Worksheets.Item($q4.Index + 1)
$newSheet.Name = "2022-Q1"

# Columns B:G hold text-like values (fund code keeps leading zeros, the
# numeric-looking figures are stored as text in this workbook) - force text
# formatting before writing so PowerShell/Excel doesn't coerce them to
# numbers.
$newSheet.Range("B2:G5").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "009778"
$newSheet.Range("C2").Value = "长信消费升级混合A"
$newSheet.Range("D2").Value = "3.43"
$newSheet.Range("E2").Value = "83.25"
$newSheet.Range("F2").Value = "5.01"
$newSheet.Range("G2").Value = "0.1718"
$newSheet.Range("H2").Value = 5

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "519959"
$newSheet.Range("C3").Value = "长信多利灵活配置混合"
$newSheet.Range("D3").Value = "1.45"
$newSheet.Range("E3").Value = "85.11"
$newSheet.Range("F3").Value = "4.75"
$newSheet.Range("G3").Value = "0.0689"
$newSheet.Range("H3").Value = 4

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "013488"
$newSheet.Range("C4").Value = "长信多利灵活配置混合D"
$newSheet.Range("D4").Value = "1.45"
$newSheet.Range("E4").Value = "85.11"
$newSheet.Range("F4").Value = "4.75"
$newSheet.Range("G4").Value = "0.0689"
$newSheet.Range("H4").Value = 4

$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "519987"
$newSheet.Range("C5").Value = "长信恒利优势混合"
$newSheet.Range("D5").Value = "0.22"
$newSheet.Range("E5").Value = "82.39"
$newSheet.Range("F5").Value = "4.49"
$newSheet.Range("G5").Value = "0.0099"
$newSheet.Range("H5").Value = 7

# ---------------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: push every existing data row down by
#    one and put the new 2022-Q1 figures in the now-empty row 2. Column A's
#    running index (0,1,2,3,...) already lines up row-for-row with the new
#    layout, so only the brand-new last row (6) needs a fresh index cell.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Make room: row 6 doesn't exist yet, clone row 5's column-A formatting into
# it before filling values so it gets the same bold/centered/bordered style.
$total.Range("A5").Copy()
$total.Range("A6").PasteSpecial(-4122)

$total.Range("B6").Value = "2021-Q1"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 0.01
$total.Range("A6").Value = 4

$total.Range("B5").Value = "2021-Q2"
$total.Range("C5").Value = 7
$total.Range("D5").Value = 0.41

$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 8
$total.Range("D4").Value = 2.46

$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 4
$total.Range("D3").Value = 0.62

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.32

# ---------------------------------------------------------------------------
# 3) Restore the originally active sheet/selection so tab state is unchanged.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
$wb.Worksheets.Item("2021-Q1").Range("A1").Select()
